$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 158.370640039444

$ws.Range("A4").Value = 18814.618
$ws.Range("B4").Value = 18311

$ws.Range("F4").Value = 6411.518
$ws.Range("G4").Value = 6359
